$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell F1 "time_taken", copying the formatting from E1 (bold/border/centered header style)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F1").Value = "time_taken"

# Populate the time_taken values for each data row (plain text, no special formatting)
$times = @(
    "2021-10-05 10:51:12.309102",
    "2021-10-05 10:51:12.309112",
    "2021-10-05 10:51:12.309115",
    "2021-10-05 10:51:12.309118",
    "2021-10-05 10:51:12.309120",
    "2021-10-05 10:51:12.309123",
    "2021-10-05 10:51:12.309126",
    "2021-10-05 10:51:12.309129",
    "2021-10-05 10:51:12.309131",
    "2021-10-05 10:51:12.309134",
    "2021-10-05 10:51:12.309137",
    "2021-10-05 10:51:12.309139",
    "2021-10-05 10:51:12.309142",
    "2021-10-05 10:51:12.309144",
    "2021-10-05 10:51:12.309147",
    "2021-10-05 10:51:12.309149",
    "2021-10-05 10:51:12.309152"
)

for ($i = 0; $i -lt $times.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $times[$i]
}
